# Auto-generated edit script for Fonds de solidarite volet 2 - add 2020-08-04 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# 1) In-place updates to nombre_aides (C) / montant_total (D) for existing region/category rows
$inplaceUpdates = @(
    ,@(28, "230", "585548.00")
    ,@(30, "468", "1746312.70")
    ,@(32, "343", "1087354.96")
    ,@(45, "287", "779606.74")
    ,@(47, "532", "1928644.11")
    ,@(48, "347", "1121385.16")
    ,@(51, "2783", "6175500.33")
    ,@(52, "20", "118500.00")
    ,@(53, "3455", "10797539.50")
    ,@(55, "3587", "10413274.64")
    ,@(56, "52", "139350.00")
    ,@(57, "67", "214868.00")
)
foreach ($u in $inplaceUpdates) {
    $r = $u[0]
    Set-TextCell $ws.Cells.Item($r, 3) $u[1]
    Set-TextCell $ws.Cells.Item($r, 4) $u[2]
}

# 2) Insert 3 new rows before row 58 for "La Réunion" (reg 04), shifting subsequent rows down
$ws.Range("A58:H60").EntireRow.Insert()

# 3) Populate the newly inserted rows
$newRows = @(
    ,@(58, "Fonds de solidarité", "VOLET2", "22", "56000.00", "04", "La Réunion", "10", "Entrepreneur individuel")
    ,@(59, "Fonds de solidarité", "VOLET2", "50", "162027.00", "04", "La Réunion", "54", "Société à responsabilité limitée (SARL)")
    ,@(60, "Fonds de solidarité", "VOLET2", "10", "28277.00", "04", "La Réunion", "57", "Société par actions simplifiée")
)
foreach ($nr in $newRows) {
    $r = $nr[0]
    for ($col = 1; $col -le 8; $col++) {
        Set-TextCell $ws.Cells.Item($r, $col) $nr[$col]
    }
}

# 4) Additional nombre_aides (C) / montant_total (D) updates on rows shifted down by the insert
$postShiftUpdates = @(
    ,@(65, "7", "31900.00")
    ,@(66, "5", "13000.00")
    ,@(72, "347", "848635.70")
    ,@(74, "846", "2784581.34")
    ,@(75, "483", "1498920.87")
    ,@(76, "33", "88500.00")
    ,@(84, "197", "459071.00")
    ,@(86, "466", "1533862.50")
    ,@(87, "173", "485976.09")
)
foreach ($u in $postShiftUpdates) {
    $r = $u[0]
    Set-TextCell $ws.Cells.Item($r, 3) $u[1]
    Set-TextCell $ws.Cells.Item($r, 4) $u[2]
}

Write-Host "Done applying Fonds de solidarite 2020-08-04 update."
